# Header cell D1 currently holds a rich-text string "LCSC Part #（optional）"
# made of four differently-formatted runs (Arial "LCSC Part #", SimSun "（",
# Arial "optional", SimSun "）"). Simplify it down to plain "LCSC Part #",
# keeping the cell's existing (bold Arial) formatting/style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "LCSC Part #"

# Move/restore the saved selection to D1 (previously it was parked on A14).
$ws.Range("D1").Select()
